$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Sheet name:" $ws.Name
